$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (table row 1): 92×54=, 26×54=, 50×66=, 47×12=, 52×60=
$t.Cell(1, 1).Range.Text = "70×68="
$t.Cell(1, 2).Range.Text = "21×43="
$t.Cell(1, 3).Range.Text = "36×43="
$t.Cell(1, 4).Range.Text = "48×61="
$t.Cell(1, 5).Range.Text = "25×46="

# Row 5 (table row 5): 13×74=, 95×45=, 82×48=, 69×61=, 43×19=
$t.Cell(5, 1).Range.Text = "88×47="
$t.Cell(5, 2).Range.Text = "32×84="
$t.Cell(5, 3).Range.Text = "96×74="
$t.Cell(5, 4).Range.Text = "27×49="
$t.Cell(5, 5).Range.Text = "48×40="

# Row 10 (table row 10): 66×49=, 93×59=, 62×26=, 57×25=, 31×85=
$t.Cell(10, 1).Range.Text = "93×88="
$t.Cell(10, 2).Range.Text = "68×60="
$t.Cell(10, 3).Range.Text = "87×62="
$t.Cell(10, 4).Range.Text = "22×43="
$t.Cell(10, 5).Range.Text = "74×59="

# Row 15 (table row 15): 99×35=, 31×85=, 90×69=, 86×55=, 26×75=
$t.Cell(15, 1).Range.Text = "82×52="
$t.Cell(15, 2).Range.Text = "88×35="
$t.Cell(15, 3).Range.Text = "94×80="
$t.Cell(15, 4).Range.Text = "39×25="
$t.Cell(15, 5).Range.Text = "26×95="

# Row 20 (table row 20): 99×65=, 61×66=, 65×37=, 93×12=, 74×36=
$t.Cell(20, 1).Range.Text = "46×52="
$t.Cell(20, 2).Range.Text = "66×16="
$t.Cell(20, 3).Range.Text = "11×14="
$t.Cell(20, 4).Range.Text = "35×46="
$t.Cell(20, 5).Range.Text = "79×66="
